$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "U2IK1P"
$ws.Range("B18").Value = "Cabezal Epson"
$ws.Range("C18").Value = "TM U950"
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 280000
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 6
$ws.Range("H18").Formula = "=(E18-D18)*G18"
$ws.Range("I18").Formula = "=D18*F18"
$ws.Range("J18").Value = 0
